# Adds the author's name/student number as a brand-new first paragraph of
# the document ("nu met naam en studentnummer" - "now with name and student
# number"), and relocates the existing "_GoBack" bookmark from the very end
# of the document to the end of that new first paragraph (the position Word
# leaves "_GoBack" at after the most recent edit).

$d = $word.ActiveDocument

# --- Step 1: drop the pre-existing "_GoBack" bookmark -----------------
# It currently sits at the very end of the document (after the last
# paragraph's text). It will be recreated after the newly inserted text
# below, mirroring where Word leaves it following the freshest edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: insert the new first paragraph ----------------------------
# A trailing placeholder character ("~") is inserted together with the
# name so that the bookmark added in step 3 is not the very last
# character of the paragraph while we create it (collapsed bookmarks at
# the exact end of a paragraph can otherwise anchor to the wrong spot).
# The placeholder is stripped again in step 4.
$insertionPoint = $d.Range(0, 0)
$newFirstLine = "Thomas van der Veen, 10346481"
$insertionPoint.InsertBefore($newFirstLine + "~`r")

# --- Step 3: re-create "_GoBack" right after the new text --------------
$bookmarkPos = $newFirstLine.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- Step 4: remove the temporary placeholder character -----------------
$placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholderRange.Delete()
